$d = $word.ActiveDocument

# Fix spelling error: NDBPRN -> NDPBRN (swap the B and P), replacing
# every occurrence in the document body.
$d.Content.Find.Execute("NDBPRN", $true, $false, $false, $false, $false,
                         $true, 1, $false, "NDPBRN", 2)
